# ProjectExpenses.xlsx — "Working single department schedule generation"
#
# 1. Mark the workbook window as minimized (best-effort; the host engine
#    does not currently round-trip bookViews/@minimized, but WindowState
#    is the closest semantic analog so we still set it).
# 2. Fill in 4 previously-blank expense rows (82-85) with real data,
#    which also introduces 4 new shared strings and grows the wrapped
#    E-column row heights.
# 3. Push the old blank rows + totals row (formerly row 88) down by 6
#    rows so the sheet keeps 6 blank rows before the total, now on row 94.
# 4. Update the total formula so it sums through the new last data row.
# 5. Leave the selection/viewport where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
    $excel.ActiveWindow.WindowState = -4140   # xlMinimized
} catch {
}

# Push rows 88.. down by 6 so the old "blank row" block (82-87) plus the
# totals row (88) end up at 88-93 (blank) and 94 (totals) respectively.
$ws.Rows.Item(88).Resize(6).Insert() | Out-Null

# Row 82 — Внедряване и доработки на място
$ws.Range("A82").Value = 42450
$ws.Range("C82").Value = 4
$ws.Range("D82").Value = 200

# Row 83 — Корекции по функционалността за редакция на екипи
$ws.Range("A83").Value = 42451
$ws.Range("C83").Value = 2
$ws.Range("D83").Value = 100
$ws.Rows.Item(83).RowHeight = 30

# Row 84 — Корекции по различни функционалности...
$ws.Range("A84").Value = 42452
$ws.Range("C84").Value = 4
$ws.Range("D84").Value = 200
$ws.Rows.Item(84).RowHeight = 45

# Row 85 — Разработка на потребителски интерфейс...
$ws.Range("A85").Value = 42453
$ws.Range("C85").Value = 6
$ws.Range("D85").Value = 300
$ws.Rows.Item(85).RowHeight = 75

# Write the new shared strings in the same order the author introduced
# them (99: row82, 100: row85, 101: row83, 102: row84) so the
# sharedStrings table comes out in the same order as the target.
$ws.Range("E82").Value = "Внедряване и доработки на място"
$ws.Range("E85").Value = "Разработка на потребителски интерфейс за отразяване на личен график и присъствена форма. Отваряне на досието от прозореца графици и оцветяване на колоните за неработни дни."
$ws.Range("E83").Value = "Корекции по функционалността за редакция на екипи"
$ws.Range("E84").Value = "Корекции по различни функционалности, задаване на работно време и списък от слъжности за вяско от сменните звена."

# Totals row moved from 88 -> 94; extend the SUM to cover the new rows.
$ws.Range("D94").Formula = "=SUM(D2:D85)"

# Leave the viewport/selection where the author was last working.
$ws.Range("A79").Select() | Out-Null
$ws.Range("E85").Select() | Out-Null
